$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Table indent: tblInd w="-5" -> w="-10" (points: -0.25 -> -0.5) ---
$t.Rows.LeftIndent = -0.5

# --- Table-level default cell margin (tblCellMar left): 103 -> 98 dxa (5.15 -> 4.9 pt) ---
$t.LeftPadding = 4.9

# --- Column widths (gridCol / tcW), in points (dxa / 20) ---
# Before: 1460, 1461, 1457, 1461, 1461, 2197
# After : 1459, 1461, 1457, 1460, 1461, 2198
$t.Cell(1,1).Width = 72.95
$t.Cell(1,4).Width = 73
$t.Cell(1,6).Width = 109.9

# --- Per-cell left margins (tcMar left): 103 -> 98 dxa (5.15 -> 4.9 pt) for every cell ---
for ($r = 1; $r -le $t.Rows.Count; $r++) {
  for ($c = 1; $c -le $t.Columns.Count; $c++) {
    $t.Cell($r, $c).LeftPadding = 4.9
  }
}

# --- Updated perf numbers in the "Lecture d'un fichier (LOAD)" row (row 3) ---
$row3 = $t.Rows.Item(3)

$rng = $t.Cell(3,2).Range
$rng.Find.Execute("1,311", $true, $false, $false, $false, $false, $true, 1, $false, "1,14", 2)

$rng = $t.Cell(3,3).Range
$rng.Find.Execute("12,006", $true, $false, $false, $false, $false, $true, 1, $false, "9,383", 2)

$rng = $t.Cell(3,4).Range
$rng.Find.Execute("78,885", $true, $false, $false, $false, $false, $true, 1, $false, "69,332", 2)

$rng = $t.Cell(3,5).Range
$rng.Find.Execute("761,14", $true, $false, $false, $false, $false, $true, 1, $false, "498,009", 2)

$rng = $t.Cell(3,6).Range
$rng.Find.Execute("8347,21", $true, $false, $false, $false, $false, $true, 1, $false, "5299,82", 2)
